$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Knight's play 1-9")
$ws2 = $wb.Worksheets.Item("Knight's play 10-18")
$ws3 = $wb.Worksheets.Item("Knight's play 19-27")

# ---------------------------------------------------------------------------
# Sheet1 ("Knight's play 1-9"): scroll position only
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("A14").Select()
$excel.ActiveWindow.ScrollRow = 14

# ---------------------------------------------------------------------------
# Sheet2 ("Knight's play 10-18"): fix row 55, add rows 59-64
# ---------------------------------------------------------------------------
$ws2.Range("M55").ClearContents()
$ws2.Range("L55").Value = "W"

# Copy formats from the last full block (rows 52-57) down into the new block (59-64)
$ws2.Range("A52:M52").Copy()
$ws2.Range("A59:M59").PasteSpecial(-4122)
$ws2.Range("A53:M53").Copy()
$ws2.Range("A60:M60").PasteSpecial(-4122)
$ws2.Range("A54:M54").Copy()
$ws2.Range("A61:M61").PasteSpecial(-4122)
$ws2.Range("A55:M55").Copy()
$ws2.Range("A62:M62").PasteSpecial(-4122)
$ws2.Range("A56").Copy()
$ws2.Range("A63").PasteSpecial(-4122)
$ws2.Range("D56:L56").Copy()
$ws2.Range("D63:L63").PasteSpecial(-4122)
$ws2.Range("A57").Copy()
$ws2.Range("A64").PasteSpecial(-4122)
$ws2.Range("D57:M57").Copy()
$ws2.Range("D64:M64").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 59 - date / start time / end time
$ws2.Range("A59").Value = 44418
$ws2.Range("B59").Value = 0.5625
$ws2.Range("C59").Value = 0.63541666666666663

# Row 60 - Black tee scores
$ws2.Range("A60").Value = "Black"
$ws2.Range("D60").Value = 2
$ws2.Range("E60").Value = 3
$ws2.Range("F60").Value = 5
$ws2.Range("G60").Value = 2
$ws2.Range("H60").Value = 4
$ws2.Range("I60").Value = 4
$ws2.Range("J60").Value = 4
$ws2.Range("K60").Value = 5
$ws2.Range("L60").Value = 6
$ws2.Range("M60").Formula = "=SUM(D60:L60)"

# Row 61 - putts
$ws2.Range("A61").Value = "putts"
$ws2.Range("D61").Value = 1
$ws2.Range("E61").Value = 2
$ws2.Range("F61").Value = 3
$ws2.Range("G61").Value = 1
$ws2.Range("H61").Value = 2
$ws2.Range("I61").Value = 2
$ws2.Range("J61").Value = 1
$ws2.Range("K61").Value = 2
$ws2.Range("L61").Value = 3
$ws2.Range("M61").Formula = "=SUM(D61:L61)"

# Row 62 - penalties
$ws2.Range("A62").Value = "penalties"
$ws2.Range("J62").Value = "MM"
$ws2.Range("K62").Value = "W"
$ws2.Range("L62").Value = "W"

# Row 63 - Fairways Hit
$ws2.Range("A63").Value = "Fairways Hit"
$ws2.Range("D63:L63").Value = "NA"

# Row 64 - Strokes inside 80
$ws2.Range("A64").Value = "Strokes inside 80"
$ws2.Range("D64").Value = 1
$ws2.Range("E64").Value = 2
$ws2.Range("F64").Value = 4
$ws2.Range("G64").Value = 1
$ws2.Range("H64").Value = 3
$ws2.Range("I64").Value = 3
$ws2.Range("J64").Value = 3
$ws2.Range("K64").Value = 4
$ws2.Range("L64").Value = 5
$ws2.Range("L64").Font.Italic = $true
$ws2.Range("M64").Formula = "=SUM(D64:L64)"

# ---------------------------------------------------------------------------
# Sheet3 ("Knight's play 19-27"): add rows 15-20
# ---------------------------------------------------------------------------
$ws2.Range("A52:M52").Copy()
$ws3.Range("A15:M15").PasteSpecial(-4122)
$ws2.Range("A53:M53").Copy()
$ws3.Range("A16:M16").PasteSpecial(-4122)
$ws2.Range("A54:M54").Copy()
$ws3.Range("A17:M17").PasteSpecial(-4122)
$ws2.Range("A55:M55").Copy()
$ws3.Range("A18:M18").PasteSpecial(-4122)
$ws2.Range("A56").Copy()
$ws3.Range("A19").PasteSpecial(-4122)
$ws2.Range("D56:L56").Copy()
$ws3.Range("D19:L19").PasteSpecial(-4122)
$ws2.Range("A57").Copy()
$ws3.Range("A20").PasteSpecial(-4122)
$ws2.Range("D57:M57").Copy()
$ws3.Range("D20:M20").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 15 - date / start time / end time
$ws3.Range("A15").Value = 44418
$ws3.Range("B15").Value = 0.51527777777777783
$ws3.Range("C15").Value = 0.5625

# Row 16 - Black tee scores
$ws3.Range("A16").Value = "Black"
$ws3.Range("D16").Value = 5
$ws3.Range("E16").Value = 5
$ws3.Range("F16").Value = 4
$ws3.Range("G16").Value = 4
$ws3.Range("H16").Value = 4
$ws3.Range("I16").Value = 4
$ws3.Range("J16").Value = 5
$ws3.Range("K16").Value = 4
$ws3.Range("L16").Value = 4
$ws3.Range("M16").Formula = "=SUM(D16:L16)"

# Row 17 - putts
$ws3.Range("A17").Value = "putts"
$ws3.Range("D17").Value = 2
$ws3.Range("E17").Value = 2
$ws3.Range("F17").Value = 2
$ws3.Range("G17").Value = 2
$ws3.Range("H17").Value = 2
$ws3.Range("I17").Value = 2
$ws3.Range("J17").Value = 2
$ws3.Range("K17").Value = 2
$ws3.Range("L17").Value = 2
$ws3.Range("M17").Formula = "=SUM(D17:L17)"

# Row 18 - penalties
$ws3.Range("A18").Value = "penalties"
$ws3.Range("I18").Value = "M"

# Row 19 - Fairways Hit
$ws3.Range("A19").Value = "Fairways Hit"
$ws3.Range("D19:L19").Value = "NA"

# Row 20 - Strokes inside 80
$ws3.Range("A20").Value = "Strokes inside 80"
$ws3.Range("D20").Value = 4
$ws3.Range("E20").Value = 4
$ws3.Range("F20").Value = 3
$ws3.Range("G20").Value = 3
$ws3.Range("H20").Value = 3
$ws3.Range("I20").Value = 3
$ws3.Range("J20").Value = 4
$ws3.Range("K20").Value = 3
$ws3.Range("L20").Value = 3
$ws3.Range("L20").Font.Italic = $true
$ws3.Range("M20").Formula = "=SUM(D20:L20)"

# ---------------------------------------------------------------------------
# Sheet views / selections
# ---------------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("A58:M64").Select()
$excel.ActiveWindow.ScrollRow = 24

$ws3.Activate()
$ws3.Range("M20").Select()

$ws1.Range("A54").Select()

$excel.ActiveWindow.WindowState = -4143
